$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the header note in C1 with rich text: bold text, with a
# --- superscripted "th" in "5th" (was "3rd"), e.g. "Last updated 5th September 2025."
$note = $ws.Range("C1")
$note.Value = "Uzebox Omega v1.1.1 Digi-Key BOM. Last updated 5th September 2025."

# "Uzebox Omega v1.1.1 Digi-Key BOM. Last updated 5"  (48 chars, 1-based start 1)
$notePart1 = $note.Characters(1, 48)
$notePart1.Font.Bold = $true
$notePart1.Font.Name = "Tahoma"
$notePart1.Font.Size = 10
$notePart1.Font.Color = 0

# Register the bold+superscript Tahoma font (matches the "th" run's format)
# via a scratch cell so it lands in the workbook's font table, then drop the
# scratch cell so it leaves no trace in the sheet.
$scratch = $ws.Range("Z1")
$scratch.Value = "x"
$scratch.Font.Bold = $true
$scratch.Font.Superscript = $true
$scratch.Font.Name = "Tahoma"
$scratch.Font.Size = 10
$scratch.Font.Color = 0
$scratch.Clear()

# "th" (2 chars, 1-based start 49) -- superscript
$notePart2 = $note.Characters(49, 2)
$notePart2.Font.Bold = $true
$notePart2.Font.Superscript = $true
$notePart2.Font.Name = "Tahoma"
$notePart2.Font.Size = 10
$notePart2.Font.Color = 0

# " September 2025." (16 chars, 1-based start 51)
$notePart3 = $note.Characters(51, 16)
$notePart3.Font.Bold = $true
$notePart3.Font.Name = "Tahoma"
$notePart3.Font.Size = 10
$notePart3.Font.Color = 0

# --- Add the two new DigiKey BOM line items (SPI RAM + socket) at the end
# --- of the table.
$ws.Range("A26").Value = "AE9986-ND"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "CONN IC DIP SOCKET 8POS TIN        (U6)   - optional but recommended"

$ws.Range("A27").Value = "23LC1024-I/P-ND"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "IC SRAM 1MBIT SPI/QUAD 8DIP         (U6)   - optional but recommended"

# --- Reset the active selection back to A1.
$ws.Range("A1").Select()
